# mainpage improvements, error prevention
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the last logged day: it was actually 2.5h -> 3h, and its date
# was mis-entered (Jan 20th instead of Jan 24th).
$ws.Range("A12").Value = [DateTime]"2026-01-24"
$ws.Range("B12").Value = 3

# Add the missing entry for the following day (error prevention: make sure
# every worked day gets logged). Copy the date formatting from the row
# above so the new date cell matches the rest of the column.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A13").Value = [DateTime]"2026-01-25"
$ws.Range("B13").Value = 0.5

# Move the selection to the next empty row, ready for the next entry.
$ws.Range("B14").Select()
